$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 154 (1-indexed), shifting rows 154:228 down to 155:229
$ws.Rows.Item(154).Insert()

# Populate the new row 154 with the data
$ws.Cells.Item(154, 1).Value = 3
$ws.Cells.Item(154, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(154, 3).Value = "Coquimbo"
$ws.Cells.Item(154, 4).Value = 44452
$ws.Cells.Item(154, 5).Value = 5
$ws.Cells.Item(154, 6).Value = 100112017
$ws.Cells.Item(154, 7).Value = "Apio"
$ws.Cells.Item(154, 8).Value = "Americana (o)"
$ws.Cells.Item(154, 9).Value = "Primera"
$ws.Cells.Item(154, 10).Value = 260
$ws.Cells.Item(154, 11).Value = 9500
$ws.Cells.Item(154, 12).Value = 10000
$ws.Cells.Item(154, 13).Value = 9769
$ws.Cells.Item(154, 14).Value = "`$/docena de matas"
$ws.Cells.Item(154, 15).Value = "Pan de Az" + [char]0xFA + "car"
$ws.Cells.Item(154, 16).Value = 1628
$ws.Cells.Item(154, 17).Value = 6
$ws.Cells.Item(154, 18).Value = "Hortaliza"
